$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$cell = $ws.Range("A2")
$cell.Value = $cell.Value2 -replace "I prefer quit restaurants", "I prefer quiet restaurants"
